# New weekly price record for "Hortaliza, Terminal La Palmera de La Serena - Cilantro".
# A new row is inserted at row 11 (pushing the existing rows 11-127 down to 12-128),
# and the new row is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = (Get-Date -Year 2022 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112040
$ws.Cells.Item(11, 7).Value = "Cilantro"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 2400
$ws.Cells.Item(11, 11).Value = 2500
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = 2750
$ws.Cells.Item(11, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 1833
$ws.Cells.Item(11, 17).Value = 1.5
$ws.Cells.Item(11, 18).Value = "Hortaliza"
